$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Year of Treatment" column (column B) entirely, shifting
# columns C:F (All, Males, Females, Not known / missing) left to B:E.
$ws.Columns.Item(2).Delete()

# Update the remaining header labels, appending the ".deja.deja.deja" suffix.
$ws.Range("B1").Value = "All.deja.deja.deja"
$ws.Range("C1").Value = "Males.deja.deja.deja"
$ws.Range("D1").Value = "Females.deja.deja.deja"
$ws.Range("E1").Value = "Not known / missing.deja.deja.deja"
